# Insert a new column before column K ("Hyperscaler" sub-product column),
# which shifts the existing K:U columns to L:V, then populate the new
# header cell K3 with the new "Hyperscaler" label and restore the
# top-row selection to span the (now one-column-wider) merged title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns K:U one column to the right by inserting a blank column at K.
$ws.Columns.Item(11).Insert() | Out-Null

# New header text for the inserted column.
$ws.Range("K3").Value = "Hyperscaler"

# Match the column's best-fit width for the new "Hyperscaler" header text.
$ws.Columns.Item(11).ColumnWidth = 9.666666666666666

# Keep the selection in sync with the widened merged header (A1:R1 -> A1:S1).
$ws.Range("A1:S1").Select() | Out-Null
